$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "304.88"
Set-TextValue $ws.Range("E2") "1.34%"

# Row 3
Set-TextValue $ws.Range("D3") "35.80"
Set-TextValue $ws.Range("E3") "1.79%"

# Row 4
Set-TextValue $ws.Range("D4") "5.101"
Set-TextValue $ws.Range("E4") "1.36%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08076"
Set-TextValue $ws.Range("E5") "1.45%"

# Row 6
Set-TextValue $ws.Range("E6") "0.69%"

# Row 7
Set-TextValue $ws.Range("D7") "7.747"
Set-TextValue $ws.Range("E7") "-0.58%"

# Row 8
Set-TextValue $ws.Range("D8") "0.9270"
Set-TextValue $ws.Range("E8") "0.67%"

# Row 9
Set-TextValue $ws.Range("D9") "0.1350"
Set-TextValue $ws.Range("E9") "-0.31%"

# Row 10
Set-TextValue $ws.Range("D10") "0.1899"
Set-TextValue $ws.Range("E10") "2.97%"

# Row 11
Set-TextValue $ws.Range("D11") "0.09181"
Set-TextValue $ws.Range("E11") "-4.02%"

# Row 12
Set-TextValue $ws.Range("D12") "0.03418"
Set-TextValue $ws.Range("E12") "-5.50%"

# Row 13
Set-TextValue $ws.Range("D13") "0.09821"
Set-TextValue $ws.Range("E13") "-0.39%"

# Row 14
Set-TextValue $ws.Range("D14") "0.001414"
Set-TextValue $ws.Range("E14") "1.78%"

# Row 15
Set-TextValue $ws.Range("D15") "0.005910"
Set-TextValue $ws.Range("E15") "2.57%"

# Row 16
Set-TextValue $ws.Range("D16") "3.554"
Set-TextValue $ws.Range("E16") "1.24%"

# Row 17
Set-TextValue $ws.Range("D17") "4.182"
Set-TextValue $ws.Range("E17") "3.58%"

# Row 18
Set-TextValue $ws.Range("D18") "2.979"
Set-TextValue $ws.Range("E18") "1.70%"

# Row 19
Set-TextValue $ws.Range("D19") "0.3454"
Set-TextValue $ws.Range("E19") "0.85%"

# Row 20
Set-TextValue $ws.Range("D20") "0.1333"
Set-TextValue $ws.Range("E20") "1.70%"

# Row 21
Set-TextValue $ws.Range("D21") "4.899"
Set-TextValue $ws.Range("E21") "-3.06%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2603"
Set-TextValue $ws.Range("E22") "5.54%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04457"
Set-TextValue $ws.Range("E23") "-1.05%"

# Row 24
Set-TextValue $ws.Range("E24") "0.39%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004803"
Set-TextValue $ws.Range("E25") "0.17%"

# Row 26
Set-TextValue $ws.Range("E26") "3.84%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0003134"
Set-TextValue $ws.Range("E27") "4.16%"

# Row 39
Set-TextValue $ws.Range("D39") "0.01979"
Set-TextValue $ws.Range("E39") "5.49%"

# Row 40
Set-TextValue $ws.Range("D40") "0.04911"
Set-TextValue $ws.Range("E40") "4.54%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007664"
Set-TextValue $ws.Range("E41") "1.45%"

# Row 42
Set-TextValue $ws.Range("D42") "0.009161"
Set-TextValue $ws.Range("E42") "-5.68%"

# Row 43
Set-TextValue $ws.Range("D43") "0.1376"
Set-TextValue $ws.Range("E43") "3.92%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002103"
Set-TextValue $ws.Range("E44") "-0.58%"

# Row 45
Set-TextValue $ws.Range("D45") "0.01136"
Set-TextValue $ws.Range("E45") "5.26%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006391"
Set-TextValue $ws.Range("E46") "3.01%"

# Row 47
Set-TextValue $ws.Range("E47") "-0.09%"

# Row 48
Set-TextValue $ws.Range("D48") "64.67"
Set-TextValue $ws.Range("E48") "0.29%"

# Row 49
Set-TextValue $ws.Range("D49") "0.001193"
Set-TextValue $ws.Range("E49") "-19.93%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002103"
Set-TextValue $ws.Range("E50") "-0.09%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002003"
Set-TextValue $ws.Range("E51") "-0.09%"
